$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Gemeldete Krankenkasse (vollständiger Name)"
$ws.Range("A3").Value = "Gemeldete Krankenkasse (Abkürzung)"
$ws.Range("B7").Value = "01.01.2024"

$ws.Range("A10").Select()
